$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- "Model 1" textbox: "e.g. Random Forest Classifier" -> "Logistic Regression"
$model1 = $s.Shapes.Item(18)
if ($model1.TextFrame.TextRange.Runs(1,1).Text -notlike "Model 1*") {
    throw "Shape 18 is not the 'Model 1' textbox: $($model1.TextFrame.TextRange.Text)"
}
$model1.TextFrame.TextRange.Runs(1,1).Text = "Model 1 – Logistic Regression"

# --- Remove the old "Model 2 – e.g. Random Forest Classifier" textbox entirely
$oldModel2 = $s.Shapes.Item(24)
if ($oldModel2.TextFrame.TextRange.Runs(1,1).Text -notlike "Model 2*") {
    throw "Shape 24 is not the old 'Model 2' textbox: $($oldModel2.TextFrame.TextRange.Text)"
}
$oldModel2.Delete()

# --- "Model 3" textbox (now shifted up one slot): "e.g. Random Forest Classifier" -> "Random Forest Classifier"
$model3 = $s.Shapes.Item(24)
if ($model3.TextFrame.TextRange.Runs(1,1).Text -notlike "Model 3*") {
    throw "Shape 24 is not the 'Model 3' textbox: $($model3.TextFrame.TextRange.Text)"
}
$model3.TextFrame.TextRange.Runs(1,1).Text = "Model 3 – Random Forest Classifier"

# --- Add a fresh "Model 2 – SVM w/ Gaussian Kernel" textbox at the end of the shape
# tree, built from a duplicate of the "Model 1" box so all of its formatting
# (fill/line/extLst/paragraph styling) carries over exactly.
$newModel2Range = $model1.Duplicate()
$newModel2 = $newModel2Range.Item(1)
$newModel2.Left = 13462000 / 12700.0
$newModel2.Top = 19431895 / 12700.0

$newTr = $newModel2.TextFrame.TextRange
$newTr.Runs(1,1).Text = "Model 2 – SVM w/ Gaussian Kernel"
$newTr.Runs(2,1).Text = "Motivation and performance for each"

Write-Output "Done: slide now has $($s.Shapes.Count) shapes"
